# Update "想去人数" (number of people interested) figures for two conan/漫展 events
# that are listed on both the "展览" sheet and the "全部类型" sheet.
#   南宁·熊喵M动漫嘉年华【免费】   : 1121 -> 1123
#   南宁·第二届北极光动漫展        : 2549 -> 2556

$wb = $excel.ActiveWorkbook

# Sheet "展览": rows 3 and 4, column F
$wsExhibition = $wb.Worksheets.Item("展览")
$wsExhibition.Range("F3").Value = 1123
$wsExhibition.Range("F4").Value = 2556

# Sheet "全部类型": rows 5 and 6, column F
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F5").Value = 1123
$wsAll.Range("F6").Value = 2556
